# Completed the xls xlsx and csv functionality
# Append a new "book" record (row 7) to Sheet1, mirroring the existing
# UserId / Name / DisplayOrder / Genre / ISBN / Author / Publisher columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 107
$ws.Range("B7").Value = "New Book Chapter"
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "Technology"
$ws.Range("E7").Value = "353-890-545-131"
$ws.Range("F7").Value = "Emmanuel"
$ws.Range("G7").Value = "Emma Circle"

# Matches the author's final cell selection left in the saved workbook.
$ws.Range("G9").Select() | Out-Null
